$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "upload_seq2seq"
$ws2 = $wb.Worksheets.Item(2)   # "french"

# A stray cell with two spaces, written first so it lands early in the
# shared-string table (matches the order of the uploaded workbook).
$ws2.Range("C8").Value = "  "

# --- Sheet 1: replace the example text/label pairs ---------------------
$ws1.Range("A2").Value = "Web Applications"
$ws1.Range("A3").Value = "Data algorithms."
$ws1.Range("A4").Value = "Development staging."
$ws1.Range("B2").Value = "Software Engineering"
$ws1.Range("B3").Value = "Data Science & AI"
$ws1.Range("B4").Value = "Data Science & AI"

# --- Sheet 2: same replacement -----------------------------------------
$ws2.Range("A2").Value = "Web Applications"
$ws2.Range("A3").Value = "Data algorithms."
$ws2.Range("A4").Value = "Development staging."
$ws2.Range("B2").Value = "Software Engineering"
$ws2.Range("B3").Value = "Data Science & AI"
$ws2.Range("B4").Value = "Data Science & AI"

# Page orientation explicitly set on the "french" sheet.
$ws2.PageSetup.Orientation = 1

# --- Selections / active sheet -----------------------------------------
# sheet2 ends up with A2:B4 selected (active cell A2) and is no longer the
# active tab; sheet1 becomes the active/selected tab with the default
# top-left selection.
$ws2.Range("A2:B4").Select() | Out-Null
$ws1.Select() | Out-Null
$ws1.Range("A1").Select() | Out-Null
